# Adds a new "Qualite du sequencage" (H) column and populates the full
# sequencing-quality metadata table (rows 2-14) produced by the GPU LLM
# pipeline run.

function Set-TextCell($ws, $row, $col, $val) {
    # Force text storage so numeric- and percent-looking values ("20%", "30", ...)
    # are kept as literal strings instead of being parsed into numbers/percentages.
    $ws.Cells.Item($row, $col).NumberFormat = "@"
    $ws.Cells.Item($row, $col).Value = $val
    $ws.Cells.Item($row, $col).ClearFormats()
}

function Set-TextRow($ws, $rowNum, $vals) {
    for ($c = 1; $c -le $vals.Length; $c++) {
        Set-TextCell $ws $rowNum $c $vals[$c - 1]
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column H header: "Qualite du sequencage" (with cedilla), reusing the
#     bold/bordered header style already applied to A1:G1 (copy format from G1).
$ws.Cells.Item(1, 7).Copy() | Out-Null
$ws.Cells.Item(1, 8).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 8).Value = "Qualité du séquençage"

# --- Data rows 2-14: Examen, N° du prélèvement, Panel, Origine du prélèvement,
#     Type de prélèvement, Qualité du séquencage, % cellules, Qualité du séquençage
$row2 = @("24EM03355","24CU550062-2ebus","CLP","CurePath","Carcinome non à petites cellules NOS","Optimale","20%","")
$row3 = @("24EM03456","24CU052383","CLP","Curepath","Adénocarcinome TTF1+","Optimale","10%","")
$row4 = @("24EM03461","24CU002162-4","GP","Curepath","Adénocarcinome colorectal métastatique","Optimale","20%","")
$row5 = @("24EM03462","24219576 1.1","GP","CMP","Métastase hépatique d’un adénocarcinome mammaire","Optimale","30","")
$row6 = @("24EM03839","24EC09559","OST","Erasme","PF2","","70","Optimale")
$row7 = @("24EM04099","24CU062291-frottis2","OST","CurePath","PF1","Optimale","10","")
$row8 = @("24EM04107","24CU062294-1","OST","CurePath","PF1 oncocytaire","Optimale","10%","")
$row9 = @("24EM04337","2431646-1.1","OST","CMP Pathology","masse gastrique","Optimale","70","")
$row10 = @("24EM04347","23CU032757-1.02","OST","CurePath","carcinome urothélial invasif","Optimale","50","")
$row11 = @("24EM03451","24BB11466","GP","HUB – Centre d’Anatomie Pathologique –","Tumeur de la granulosa","Optimale","30%","")
$row12 = @("24EM03460","24MH9721 BN","GP","Centre Hospitalier de Mouscron","Adénocarcinome lieberkühnien","Optimale","50","")
$row13 = @("24EM03308","24218507-1.1","GP","CMP","Adénocarcinome pulmonaire","Optimale","50%","")
$row14 = @("24EM03352","24MH9794 RF","GP","Centre Hospitalier de Mouscron","Adénocarcinome lieberkühnien","Optimale","20%","")

Set-TextRow $ws 2 $row2
Set-TextRow $ws 3 $row3
Set-TextRow $ws 4 $row4
Set-TextRow $ws 5 $row5
Set-TextRow $ws 6 $row6
Set-TextRow $ws 7 $row7
Set-TextRow $ws 8 $row8
Set-TextRow $ws 9 $row9
Set-TextRow $ws 10 $row10
Set-TextRow $ws 11 $row11
Set-TextRow $ws 12 $row12
Set-TextRow $ws 13 $row13
Set-TextRow $ws 14 $row14
